$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list - refresh Price (D) and Volume(1h) (E) columns
$ws.Range("D2").Value = "66.813.77"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").Value = "3.098.27"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "579.73"
$ws.Range("E5").Value = "  +1.46%  "
$ws.Range("D6").Value = "172.69"
$ws.Range("E6").Value = "  +5.34%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.093.66"
$ws.Range("E8").Value = "  +4.96%  "
$ws.Range("E9").Value = "  +1.28%  "
$ws.Range("D10").Value = "6.42"
$ws.Range("E10").Value = "  -3.43%  "
$ws.Range("E11").Value = "  +2.83%  "
$ws.Range("E12").Value = "  +3.41%  "
$ws.Range("E13").Value = "  +1.99%  "
$ws.Range("D14").Value = "37.46"
$ws.Range("E14").Value = "  +6.09%  "
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").Value = "3.608.40"
$ws.Range("E16").Value = "  +4.84%  "
$ws.Range("D17").Value = "66.776.58"
$ws.Range("E17").Value = "  +1.43%  "
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("D19").Value = "3.097.37"
$ws.Range("E19").Value = "  +5.03%  "
$ws.Range("D20").Value = "16.24"
$ws.Range("E20").Value = "  +1.69%  "
$ws.Range("D21").Value = "480.46"
$ws.Range("E21").Value = "  +7.46%  "
$ws.Range("E22").Value = "  +2.62%  "
$ws.Range("E23").Value = "  +3.57%  "
$ws.Range("D24").Value = "84.03"
$ws.Range("D25").Value = "13.28"
$ws.Range("E25").Value = "  +7.63%  "
$ws.Range("E26").Value = "  +5.39%  "
$ws.Range("D27").Value = "10.05"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("E29").Value = "  -3.05%  "
$ws.Range("E30").Value = "  -1.80%  "
$ws.Range("E31").Value = "  +3.53%  "
$ws.Range("E32").Value = "  +5.41%  "
$ws.Range("D33").Value = "0.0000101"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("E34").Value = "  -2.19%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  +2.92%  "
$ws.Range("D37").Value = "0.992"
$ws.Range("E37").Value = "  +2.27%  "
$ws.Range("D38").Value = "48.05"
$ws.Range("E38").Value = "  +1.21%  "
$ws.Range("D39").Value = "2.12"
$ws.Range("E39").Value = "  +7.05%  "
$ws.Range("E40").Value = "  +4.98%  "
$ws.Range("D41").Value = "50.07"
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("D42").Value = "0.122"
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("E43").Value = "  +2.34%  "
$ws.Range("D44").Value = "2.81"
$ws.Range("E44").Value = "  -0.44%  "
$ws.Range("D45").Value = "2.841.74"
$ws.Range("E45").Value = "  +6.10%  "
$ws.Range("E46").Value = "  +2.63%  "
$ws.Range("D47").Value = "384.16"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").Value = "135.56"
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").Value = "25.13"
$ws.Range("E50").Value = "  +5.01%  "
$ws.Range("E51").Value = "  +2.63%  "
